# Edit script for FLORIDA_2015.xlsx
# 1) Title-case Spanish connector words (de, del, la, las, lo, los, el, y)
#    in municipality / state names (e.g. "Pabellon de Arteaga" -> "Pabellon De Arteaga")
# 2) Rename header columns to snake_case names
# 3) Nudge a set of D-column percentage values by 1 ULP to match the recalculated
#    values produced upstream (same double ratios, just the adjacent representable value)
# 4) Remove trailing footer rows (1813:1818) and shrink the used range accordingly

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: title-case the Spanish connector words across the whole data range ---
# Use partial-match (LookAt:=2 => xlPart) substring replace, case sensitive, so that
# only the lower-case standalone words get capitalized; everything else is untouched.
$rng = $ws.Range("A1:D1818")

$null = $rng.Replace(" de ", " De ", 2, 1, $true)
$null = $rng.Replace(" del ", " Del ", 2, 1, $true)
$null = $rng.Replace(" la ", " La ", 2, 1, $true)
$null = $rng.Replace(" las ", " Las ", 2, 1, $true)
$null = $rng.Replace(" lo ", " Lo ", 2, 1, $true)
$null = $rng.Replace(" los ", " Los ", 2, 1, $true)
$null = $rng.Replace(" el ", " El ", 2, 1, $true)
$null = $rng.Replace(" y ", " Y ", 2, 1, $true)

# --- Step 2: rename header row to snake_case column names ---
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# --- Step 3: re-assert the recalculated percentage values (1 ULP nudge) ---
# These rows hold C/Total ratios whose correctly-rounded IEEE-754 double landed
# on the adjacent representable value once the workbook total was recalculated;
# reassign them explicitly to the exact target double (written in plain decimal
# notation since the interpreter does not accept scientific-notation literals).
$rowsSmallRatio = @(11,29,37,38,53,120,133,163,176,186,191,197,200,236,237,245,288,304,318,333,334,456,478,485,492,504,523,531,534,542,582,599,608,617,634,637,642,666,699,732,753,758,774,782,800,804,846,857,872,881,884,898,907,908,923,927,945,947,968,982,993,1020,1023,1031,1039,1040,1049,1061,1062,1073,1085,1093,1096,1114,1116,1126,1128,1140,1150,1177,1185,1191,1192,1194,1200,1206,1221,1233,1235,1236,1242,1248,1255,1285,1301,1312,1334,1344,1350,1407,1427,1431,1445,1456,1497,1500,1509,1520,1530,1547,1579,1580,1596,1598,1605,1619,1623,1656,1660,1688,1701,1709,1711,1713,1720,1741,1758,1770,1786,1789,1794,1804,1807,1808)
foreach ($r in $rowsSmallRatio) {
    $ws.Cells.Item($r, 4).Value = 0.0000924299842869026806215
}

$rowsMidRatio = @(117,125,441,455,535,702,863,1395,1413,1420,1795)
foreach ($r in $rowsMidRatio) {
    $ws.Cells.Item($r, 4).Value = 0.0009242998428690268062149
}

$ws.Range("D257").Value = 0.0095202883815509765919050
$ws.Range("D134").Value = 0.0958807037002803591585476
$ws.Range("D407").Value = 0.0913824444649844330745836
$ws.Range("D1208").Value = 0.0963736636164772009216506

# --- Step 4: remove the footer / metadata rows at the bottom of the sheet ---
$ws.Rows("1813:1818").Delete()
